$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7042503952980042
$ws.Range("B1").Value = 0.6946749091148376
$ws.Range("C1").Value = 0.7280676960945129
$ws.Range("D1").Value = 0.9441881775856018
$ws.Range("E1").Value = 0.8855565190315247
